$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing existing rows 3-7 down to 4-8
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the new observation record
$ws.Range("A3").Value = 130611551
$ws.Range("B3").Value = 57884
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"

$ws.Range("L3").Value = "hona"
$ws.Range("M3").Value = "födosökande"
$ws.Range("N3").Value = "observerad"

$ws.Range("P3").Value = "Västansjö, Vb"
$ws.Range("Q3").Value = 745582
$ws.Range("R3").Value = 7101809
$ws.Range("S3").Value = 25

$ws.Range("T3").Value = "Västerbotten"
$ws.Range("U3").Value = "Umeå"
$ws.Range("V3").Value = "Västerbotten"
$ws.Range("W3").Value = "Umeå socken"

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2026-01-08"
$ws.Range("Z3").Value = "12:50"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2026-01-08"
$ws.Range("AB3").Value = "12:50"
$ws.Range("AC3").Value = "En individ födosökte på gran"

$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

$ws.Range("AW3").Value = "Alva Danielsson"
$ws.Range("AX3").Value = "Alva Danielsson"
